# A new weekly price record was added to the Piña (pineapple) price log.
# It is inserted as row 119, shifting the existing rows 119-202 down to
# 120-203 (dimension grows from A1:T202 to A1:T203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 119, pushing everything below down.
$ws.Rows.Item(119).Insert()

# Populate the new row with the latest weekly observation.
$ws.Range("A119").Value = 5
$ws.Range("B119").Value = "Macroferia Regional de Talca"
$ws.Range("C119").Value = "Maule"
$ws.Range("D119").Value = 44582
$ws.Range("E119").Value = 7
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100108
$ws.Range("H119").Value = "Tropicales y subtropicales"
$ws.Range("I119").Value = 100108005
$ws.Range("J119").Value = "Pi" + [char]0x00F1 + "a"
$ws.Range("K119").Value = "Caramelo"
$ws.Range("L119").Value = "Tercera"
$ws.Range("M119").Value = 200
$ws.Range("N119").Value = 16000
$ws.Range("O119").Value = 16000
$ws.Range("P119").Value = 16000
$ws.Range("Q119").Value = "`$/caja 16 unidades"
$ws.Range("R119").Value = "Ecuador"
$ws.Range("S119").Value = 1000
$ws.Range("T119").Value = 16
